$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.585.62'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.28%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.86%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -1.32%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '336.34'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.09%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.27%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3909'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.23'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.44%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.42%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.24'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -6.35%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.824'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.37%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.822.77'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.91%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.977'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.41%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06946'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.00%  '

$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '87.66'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.20%  '

$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.17%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.04'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.26%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.580.31'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.26%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.88%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.23'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.38%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.155'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.53%  '

$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.041.23'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.42%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.77'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.09%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.42'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.01%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.094'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.54%  '

$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.023'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.08%  '

$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.32'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.72%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9691'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.18%  '

$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09354'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.45%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.362'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.44%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.479'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.47%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.347'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.68%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06156'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.87%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02202'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.40%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.170'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.694'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.19%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5706'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.95%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.14'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.37%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.436'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.51%  '

$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1792'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.72%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.244'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.77%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5374'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.96%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.72'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.59%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.07098'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.21%  '

$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.902'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.29%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '113.21'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.345'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.60%  '
